$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5942
$ws.Range("E2").Value = 1010
$ws.Range("F2").Value = 1010
$ws.Range("G2").Value = 978
$ws.Range("H2").Value = 761
$ws.Range("I2").Value = 756
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 47149
$ws.Range("L2").Value = 37911
$ws.Range("M2").Value = 9238
$ws.Range("N2").Value = 9171
$ws.Range("O2").Value = 67
$ws.Range("P2").Value = 1105
$ws.Range("Q2").Value = -1636
$ws.Range("R2").Value = -572
$ws.Range("S2").Value = 2744
$ws.Range("T2").Value = 88
$ws.Range("U2").ClearContents()
$ws.Range("V2").Value = 6650
$ws.Range("W2").Value = 16.99
$ws.Range("X2").Value = 12.8
$ws.Range("Y2").Value = 8.56
$ws.Range("Z2").Value = 1.69
$ws.Range("AA2").Value = 410.36
$ws.Range("AB2").Value = 736.07
$ws.Range("AC2").Value = 3420
$ws.Range("AD2").Value = 13.48
$ws.Range("AE2").Value = 41500
$ws.Range("AF2").Value = 1.11
$ws.Range("AG2").Value = 450
$ws.Range("AH2").Value = 0.98
$ws.Range("AI2").Value = 13.16
$ws.Range("AJ2").Value = 22099740

# Row 3
$ws.Range("D3").Value = 8233
$ws.Range("E3").Value = 2414
$ws.Range("F3").Value = 2414
$ws.Range("G3").Value = 2496
$ws.Range("H3").Value = 1900
$ws.Range("I3").Value = 1897
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 61522
$ws.Range("L3").Value = 50424
$ws.Range("M3").Value = 11098
$ws.Range("N3").Value = 11031
$ws.Range("O3").Value = 68
$ws.Range("P3").Value = 1105
$ws.Range("Q3").Value = -1526
$ws.Range("R3").Value = -458
$ws.Range("S3").Value = 1986
$ws.Range("T3").Value = 84
$ws.Range("U3").ClearContents()
$ws.Range("V3").Value = 6260
$ws.Range("W3").Value = 29.33
$ws.Range("X3").Value = 23.07
$ws.Range("Y3").Value = 18.78
$ws.Range("Z3").Value = 3.49
$ws.Range("AA3").Value = 454.33
$ws.Range("AB3").Value = 904.4
$ws.Range("AC3").Value = 8585
$ws.Range("AD3").Value = 7.22
$ws.Range("AE3").Value = 49913
$ws.Range("AF3").Value = 1.24
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 1.13
$ws.Range("AI3").Value = 8.15
$ws.Range("AJ3").Value = 22099740

# Row 4
$ws.Range("D4").Value = 9437
$ws.Range("E4").Value = 2307
$ws.Range("F4").Value = 2307
$ws.Range("G4").Value = 2359
$ws.Range("H4").Value = 1802
$ws.Range("I4").Value = 1799
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 88571
$ws.Range("L4").Value = 76090
$ws.Range("M4").Value = 12480
$ws.Range("N4").Value = 12437
$ws.Range("O4").Value = 43
$ws.Range("P4").Value = 1105
$ws.Range("Q4").Value = -4020
$ws.Range("R4").Value = -5647
$ws.Range("S4").Value = 9802
$ws.Range("T4").Value = 61
$ws.Range("U4").ClearContents()
$ws.Range("V4").Value = 12349
$ws.Range("W4").Value = 24.45
$ws.Range("X4").Value = 19.1
$ws.Range("Y4").Value = 15.33
$ws.Range("Z4").Value = 2.4
$ws.Range("AA4").Value = 609.68
$ws.Range("AB4").Value = 1029.46
$ws.Range("AC4").Value = 8142
$ws.Range("AD4").Value = 8.84
$ws.Range("AE4").Value = 56278
$ws.Range("AF4").Value = 1.28
$ws.Range("AG4").Value = 850
$ws.Range("AH4").Value = 1.18
$ws.Range("AI4").Value = 10.44
$ws.Range("AJ4").Value = 22099740

# Row 5
$ws.Range("D5").Value = 12163
$ws.Range("E5").Value = 3158
$ws.Range("F5").Value = 3158
$ws.Range("G5").Value = 3231
$ws.Range("H5").Value = 2416
$ws.Range("I5").Value = 2402
$ws.Range("J5").Value = 14
$ws.Range("K5").Value = 118566
$ws.Range("L5").Value = 103094
$ws.Range("M5").Value = 15472
$ws.Range("N5").Value = 15246
$ws.Range("O5").Value = 226
$ws.Range("P5").Value = 1105
$ws.Range("Q5").Value = -8277
$ws.Range("R5").Value = -2922
$ws.Range("S5").Value = 11440
$ws.Range("T5").Value = 85
$ws.Range("U5").ClearContents()
$ws.Range("V5").Value = 20496
$ws.Range("W5").Value = 25.97
$ws.Range("X5").Value = 19.87
$ws.Range("Y5").Value = 17.36
$ws.Range("Z5").Value = 2.32
$ws.Range("AA5").Value = 666.33
$ws.Range("AB5").Value = 1300.18
$ws.Range("AC5").Value = 10870
$ws.Range("AD5").Value = 8.06
$ws.Range("AE5").Value = 68988
$ws.Range("AF5").Value = 1.27
$ws.Range("AG5").Value = 1300
$ws.Range("AH5").Value = 1.48
$ws.Range("AI5").Value = 11.96
$ws.Range("AJ5").Value = 22099740

# Row 6
$ws.Range("D6").Value = 21467
$ws.Range("E6").Value = 2890
$ws.Range("F6").Value = 2890
$ws.Range("G6").Value = 2810
$ws.Range("H6").Value = 1932
$ws.Range("I6").Value = 1932
$ws.Range("K6").Value = 182397
$ws.Range("L6").Value = 161739
$ws.Range("M6").Value = 20658
$ws.Range("N6").Value = 20415
$ws.Range("P6").Value = 1270
$ws.Range("Q6").Value = -15001
$ws.Range("R6").Value = -7200
$ws.Range("S6").Value = 25195
$ws.Range("T6").Value = 148
$ws.Range("U6").ClearContents()
$ws.Range("V6").Value = 37990
$ws.Range("W6").Value = 13.46
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 10.84
$ws.Range("Z6").Value = 1.28
$ws.Range("AA6").Value = 782.95
$ws.Range("AB6").Value = 1527.05
$ws.Range("AC6").Value = 7750
$ws.Range("AD6").Value = 10.09
$ws.Range("AE6").Value = 80397
$ws.Range("AF6").Value = 0.97
$ws.Range("AG6").Value = 1500
$ws.Range("AH6").Value = 1.92
$ws.Range("AI6").Value = 24.69
$ws.Range("AJ6").Value = 22099740

# Row 7
$ws.Range("D7").Value = 29254
$ws.Range("E7").Value = 4430
$ws.Range("G7").Value = 4584
$ws.Range("H7").Value = 3423
$ws.Range("I7").Value = 3461
$ws.Range("K7").Value = 226709
$ws.Range("L7").Value = 208336
$ws.Range("M7").Value = 22809
$ws.Range("N7").Value = 22810
$ws.Range("P7").Value = 1269
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 15.14
$ws.Range("X7").Value = 11.7
$ws.Range("Y7").Value = 16.01
$ws.Range("Z7").Value = 1.67
$ws.Range("AA7").Value = 913.41
$ws.Range("AC7").Value = 13630
$ws.Range("AD7").Value = 5.2
$ws.Range("AE7").Value = 91633
$ws.Range("AF7").Value = 0.77
$ws.Range("AG7").Value = 1924
$ws.Range("AH7").Value = 2.71
$ws.Range("AI7").Value = 12.29

# Row 8
$ws.Range("D8").Value = 28970
$ws.Range("E8").Value = 4197
$ws.Range("G8").Value = 4274
$ws.Range("H8").Value = 3128
$ws.Range("I8").Value = 3119
$ws.Range("K8").Value = 247436
$ws.Range("L8").Value = 227554
$ws.Range("M8").Value = 25368
$ws.Range("N8").Value = 25512
$ws.Range("P8").Value = 1269
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 14.49
$ws.Range("X8").Value = 10.8
$ws.Range("Y8").Value = 12.91
$ws.Range("Z8").Value = 1.32
$ws.Range("AA8").Value = 897.03
$ws.Range("AC8").Value = 12283
$ws.Range("AD8").Value = 5.77
$ws.Range("AE8").Value = 102485
$ws.Range("AF8").Value = 0.69
$ws.Range("AG8").Value = 1933
$ws.Range("AH8").Value = 2.73
$ws.Range("AI8").Value = 13.7

# Row 9
$ws.Range("D9").Value = 30840
$ws.Range("E9").Value = 4390
$ws.Range("G9").Value = 4481
$ws.Range("H9").Value = 3286
$ws.Range("I9").Value = 3319
$ws.Range("K9").Value = 279299
$ws.Range("L9").Value = 251426
$ws.Range("M9").Value = 27870
$ws.Range("N9").Value = 27977
$ws.Range("P9").Value = 1269
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 14.24
$ws.Range("X9").Value = 10.65
$ws.Range("Y9").Value = 12.41
$ws.Range("Z9").Value = 1.25
$ws.Range("AA9").Value = 902.15
$ws.Range("AC9").Value = 13070
$ws.Range("AD9").Value = 5.42
$ws.Range("AE9").Value = 112389
$ws.Range("AF9").Value = 0.63
$ws.Range("AG9").Value = 2062
$ws.Range("AH9").Value = 2.91
$ws.Range("AI9").Value = 13.73
